$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after "Mesh Params" and rename it
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 with the tutorial "Part 2 Geo" data ---
$ws2.Range("I2").Value = "Target"
$ws2.Range("J2").Value = 300

$ws2.Range("C3").Value = "Part 2 Geo"

$ws2.Range("C4").Value = "DP ID"
$ws2.Range("D4").Value = "In Dia"
$ws2.Range("E4").Value = "In 2 Angle"
$ws2.Range("F4").Value = "Massflow"
$ws2.Range("G4").Value = "Computational Time"
$ws2.Range("H4").Value = "OutTemp Range"
$ws2.Range("I4").Value = "Average Temp Out"
$ws2.Range("J4").Value = "Error"
$ws2.Range("D4:J4").WrapText = $true
$ws2.Rows.Item(4).RowHeight = 45

$ws2.Range("C5").Value = 6
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 1500
$ws2.Range("G5").Value = 37.86
$ws2.Range("H5").Value = 0.0097046
$ws2.Range("I5").Value = 300.002
$ws2.Range("J5").Formula = "=(I5-`$J`$2)/`$J`$2"

$ws2.Range("C6").Value = 7
$ws2.Range("D6").Value = 1
$ws2.Range("E6").Value = -45
$ws2.Range("F6").Value = 1500
$ws2.Range("G6").Value = 53.687
$ws2.Range("H6").Value = 0.0072937
$ws2.Range("I6").Value = 300.079
$ws2.Range("J6").Formula = "=(I6-`$J`$2)/`$J`$2"

$ws2.Range("C7").Value = 8
$ws2.Range("D7").Value = 1.5
$ws2.Range("E7").Value = 0
$ws2.Range("F7").Value = 1500
$ws2.Range("G7").Value = 70.1
$ws2.Range("H7").Value = 0.00424194
$ws2.Range("I7").Value = 299.866
$ws2.Range("J7").Formula = "=(I7-`$J`$2)/`$J`$2"

$ws2.Range("D8").Value = 1
$ws2.Range("E8").Value = 0
$ws2.Range("F8").Value = 1600
$ws2.Range("H8").Value = 0.00564575
$ws2.Range("I8").Value = 300.003
$ws2.Range("J8").Formula = "=(I8-`$J`$2)/`$J`$2"

$ws2.Columns.Item(10).ColumnWidth = 11.140625

# --- Sheet1 loses the tab-selected state and gets a new selection ---
$ws1.Activate()
$ws1.Range("D3:K7").Select()

# --- View state: Sheet2 becomes the active / selected tab ---
$ws2.Activate()
$ws2.Range("G8").Select()
